$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire column O (header "擬調"), shifting the following
# columns (P "調整後利率", etc.) one position to the left.
$ws.Columns("O").Delete()

# Excel normally re-points the hidden "_FilterDatabase" defined name to the
# (now narrower) header range whenever the underlying range shifts; the
# automation layer doesn't do this implicitly for a plain column delete, so
# re-create the name explicitly pointing at the new A1:O1 header range while
# preserving its hidden/worksheet-scoped nature.
foreach ($existingName in $ws.Names) {
    if ($existingName.Name -like "*_FilterDatabase") {
        $existingName.Delete()
    }
}
$filterDbName = $ws.Names.Add("_xlnm._FilterDatabase", $ws.Range("A1:O1"))
$filterDbName.Visible = $false
